$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date-format hints in the header row to use slashes instead of dashes
$ws.Range("E1").Value = "Start Date (dd/mm/yyyy)"
$ws.Range("F1").Value = "Due Date (dd/mm/yyyy)"

# Match the author's final selection (cell F1)
$ws.Range("F1").Select()
